{"js": "// Replace the 25 division-fact answers in the single 5-column table.\n// The table has 20 rows: 5 \"content\" rows (each holding 5 answers) interleaved\n// with 3 blank spacer rows. We overwrite the table's full `values` grid,\n// keeping the blank rows blank, so existing run formatting (font/size) on the\n// populated cells is preserved by the host and only the text itself changes.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst blankRow = [\"\", \"\", \"\", \"\", \"\"];\n\nconst newValues = [\n  [\"82\u00f78=10, 2\", \"52\u00f79=5, 7\", \"56\u00f76=9, 2\", \"62\u00f77=8, 6\", \"26\u00f76=4, 2\"],\n  blankRow,\n  blankRow,\n  blankRow,\n  [\"78\u00f73=26, 0\", \"74\u00f75=14, 4\", \"92\u00f76=15, 2\", \"36\u00f79=4, 0\", \"94\u00f78=11, 6\"],\n  blankRow,\n  blankRow,\n  blankRow,\n  [\"52\u00f79=5, 7\", \"74\u00f79=8, 2\", \"36\u00f75=7, 1\", \"50\u00f74=12, 2\", \"18\u00f77=2, 4\"],\n  blankRow,\n  blankRow,\n  blankRow,\n  [\"83\u00f76=13, 5\", \"83\u00f77=11, 6\", \"84\u00f75=16, 4\", \"19\u00f72=9, 1\", \"92\u00f75=18, 2\"],\n  blankRow,\n  blankRow,\n  blankRow,\n  [\"48\u00f74=12, 0\", \"48\u00f77=6, 6\", \"88\u00f72=44, 0\", \"54\u00f77=7, 5\", \"25\u00f76=4, 1\"],\n  blankRow,\n  blankRow,\n  blankRow,\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the 25 division-fact answers in the single 5-column table.\n# The table has 20 rows: 5 \"content\" rows (1, 5, 9, 13, 17 in 1-based Word\n# numbering) each holding 5 answers, interleaved with blank spacer rows.\n# Assigning straight to Cell(r,c).Range.Text keeps the existing run\n# formatting (font/size) on each cell and only swaps the visible text.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n  @(1,  @(\"82\u00f78=10, 2\", \"52\u00f79=5, 7\", \"56\u00f76=9, 2\", \"62\u00f77=8, 6\", \"26\u00f76=4, 2\")),\n  @(5,  @(\"78\u00f73=26, 0\", \"74\u00f75=14, 4\", \"92\u00f76=15, 2\", \"36\u00f79=4, 0\", \"94\u00f78=11, 6\")),\n  @(9,  @(\"52\u00f79=5, 7\", \"74\u00f79=8, 2\", \"36\u00f75=7, 1\", \"50\u00f74=12, 2\", \"18\u00f77=2, 4\")),\n  @(13, @(\"83\u00f76=13, 5\", \"83\u00f77=11, 6\", \"84\u00f75=16, 4\", \"19\u00f72=9, 1\", \"92\u00f75=18, 2\")),\n  @(17, @(\"48\u00f74=12, 0\", \"48\u00f77=6, 6\", \"88\u00f72=44, 0\", \"54\u00f77=7, 5\", \"25\u00f76=4, 1\"))\n)\n\nforeach ($rowSpec in $newValues) {\n    $rowIndex = $rowSpec[0]\n    $rowValues = $rowSpec[1]\n    for ($col = 1; $col -le 5; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
